$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'320.93"
$ws.Range("E2").Value = "'6.01%"
$ws.Range("D3").Value = "'49.28"
$ws.Range("E3").Value = "'11.14%"
$ws.Range("D4").Value = "'5.319"
$ws.Range("E4").Value = "'4.27%"
$ws.Range("D5").Value = "'0.08054"
$ws.Range("E5").Value = "'4.31%"
$ws.Range("D6").Value = "'4.622"
$ws.Range("E6").Value = "'4.43%"
$ws.Range("D7").Value = "'1.337"
$ws.Range("E7").Value = "'27.53%"
$ws.Range("E8").Value = "'1.90%"
$ws.Range("D9").Value = "'0.1260"
$ws.Range("E9").Value = "'-1.13%"
$ws.Range("D10").Value = "'0.1966"
$ws.Range("E10").Value = "'5.81%"
$ws.Range("D11").Value = "'0.09608"
$ws.Range("E11").Value = "'3.26%"
$ws.Range("D12").Value = "'0.04708"
$ws.Range("E12").Value = "'13.31%"
$ws.Range("D13").Value = "'0.1047"
$ws.Range("E13").Value = "'-0.05%"
$ws.Range("D14").Value = "'0.001322"
$ws.Range("E14").Value = "'3.48%"
$ws.Range("D15").Value = "'0.04210"
$ws.Range("E15").Value = "'0.32%"
$ws.Range("D16").Value = "'0.005826"
$ws.Range("E16").Value = "'1.15%"
$ws.Range("D17").Value = "'3.343"
$ws.Range("E17").Value = "'0.07%"
$ws.Range("D18").Value = "'2.449"
$ws.Range("E18").Value = "'5.07%"
$ws.Range("D19").Value = "'0.3525"
$ws.Range("E19").Value = "'5.48%"
$ws.Range("D20").Value = "'8.034"
$ws.Range("E20").Value = "'-0.20%"
$ws.Range("D21").Value = "'0.1363"
$ws.Range("E21").Value = "'-0.46%"
$ws.Range("D22").Value = "'0.3095"
$ws.Range("E22").Value = "'-2.60%"
$ws.Range("D23").Value = "'0.001302"
$ws.Range("E23").Value = "'1.47%"
$ws.Range("D24").Value = "'0.004270"
$ws.Range("E24").Value = "'-2.99%"
$ws.Range("E25").Value = "'0.05%"
$ws.Range("D26").Value = "'0.0003540"
$ws.Range("D38").Value = "'0.02708"
$ws.Range("E38").Value = "'7.62%"
$ws.Range("D39").Value = "'0.05955"
$ws.Range("E39").Value = "'12.06%"
$ws.Range("E40").Value = "'87.38%"
$ws.Range("D41").Value = "'0.008024"
$ws.Range("E41").Value = "'3.86%"
$ws.Range("D42").Value = "'0.1466"
$ws.Range("E42").Value = "'8.46%"
$ws.Range("D43").Value = "'0.007914"
$ws.Range("E43").Value = "'7.58%"
$ws.Range("D44").Value = "'0.007885"
$ws.Range("E44").Value = "'4.83%"
$ws.Range("D45").Value = "'0.3493"
$ws.Range("E45").Value = "'15.44%"
$ws.Range("D46").Value = "'0.00006905"
$ws.Range("E46").Value = "'2.74%"
$ws.Range("E47").Value = "'0.06%"
$ws.Range("D48").Value = "'0.05949"
$ws.Range("E48").Value = "'33.44%"
$ws.Range("D49").Value = "'0.004001"
$ws.Range("E49").Value = "'-4.76%"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'0.06%"
$ws.Range("E51").Value = "'0.06%"
